# Update "paises" (countries) data sheet + provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 22:59"

# --- Update per-country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 4481813
$ws.Range("C4").Value = 48403
$ws.Range("D4").Value = 2154879
$ws.Range("E4").Value = 2174939
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 920
$ws.Range("H4").Value = 151995

# Row 21: Alemania
$ws.Range("B21").Value = 207934
$ws.Range("C21").Value = 555
$ws.Range("E21").Value = 7327

# Row 33: Suecia
$ws.Range("B33").Value = 79494
$ws.Range("C33").Value = 39
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 5702

# Row 39: Israel
$ws.Range("B39").Value = 66293
$ws.Range("C39").Value = 2308
$ws.Range("D39").Value = 32182
$ws.Range("E39").Value = 33625

# Row 48: Guatemala
$ws.Range("B48").Value = 46451
$ws.Range("C48").Value = 1142
$ws.Range("D48").Value = 33494
$ws.Range("E48").Value = 11175
$ws.Range("G48").Value = 21
$ws.Range("H48").Value = 1782

# Row 69: Camerun
$ws.Range("B69").Value = 17179
$ws.Range("C69").Value = 69
$ws.Range("E69").Value = 2249

# Rows 72/73: Chequia and Costa de Marfil swap rank (Costa de Marfil now ahead)
$ws.Range("A72").Value = "Costa de Marfil"
$ws.Range("B72").Value = 15713
$ws.Range("C72").Value = 58
$ws.Range("D72").Value = 10537
$ws.Range("E72").Value = 5078
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 98

$ws.Range("A73").Value = "Chequia"
$ws.Range("B73").Value = 15684
$ws.Range("C73").Value = 168
$ws.Range("D73").Value = 11428
$ws.Range("E73").Value = 3882
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 374

# Row 80: Estado de Palestina
$ws.Range("E80").Value = 7107
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 79

# Rows 107/108: Maldivas and Nicaragua swap rank (Nicaragua now ahead)
$ws.Range("A107").Value = "Nicaragua"
$ws.Range("B107").Value = 3672
$ws.Range("C107").Value = 233
$ws.Range("D107").Value = 2492
$ws.Range("E107").Value = 1064
$ws.Range("G107").Value = 8
$ws.Range("H107").Value = 116

$ws.Range("A108").Value = "Maldivas"
$ws.Range("B108").Value = 3506
$ws.Range("C108").Value = 137
$ws.Range("D108").Value = 2547
$ws.Range("E108").Value = 944
$ws.Range("H108").Value = 15

# Rows 110/111: Congo and Somalia swap rank (Somalia now ahead)
$ws.Range("A110").Value = "Somalia"
$ws.Range("B110").Value = 3212
$ws.Range("C110").Value = 16
$ws.Range("D110").Value = 1562
$ws.Range("E110").Value = 1557
$ws.Range("H110").Value = 93

$ws.Range("A111").Value = "Congo"
$ws.Range("B111").Value = 3200
$ws.Range("D111").Value = 829
$ws.Range("E111").Value = 2317
$ws.Range("H111").Value = 54

# Row 148: Angola
$ws.Range("B148").Value = 1000
$ws.Range("C148").Value = 50
$ws.Range("D148").Value = 266
$ws.Range("E148").Value = 687
$ws.Range("G148").Value = 6
$ws.Range("H148").Value = 47

# Row 158: Siria
$ws.Range("B158").Value = 694
$ws.Range("C158").Value = 20
$ws.Range("D158").Value = 220
$ws.Range("E158").Value = 434
